$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        $Sheet,
        [int]$Row,
        [int]$Col,
        [string]$Text
    )
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# Row 2
Set-CellText $ws 2 4 "30.296.03"
Set-CellText $ws 2 5 "  +0.06%  "

# Row 3
Set-CellText $ws 3 4 "1.868.44"
Set-CellText $ws 3 5 "  +0.21%  "

# Row 4
Set-CellText $ws 4 4 "0.9999"
Set-CellText $ws 4 5 "  -0.19%  "

# Row 5
Set-CellText $ws 5 4 "243.28"
Set-CellText $ws 5 5 "  +3.84%  "

# Row 7
Set-CellText $ws 7 4 "0.4722"
Set-CellText $ws 7 5 "  +0.60%  "

# Row 8
Set-CellText $ws 8 2 "Cardano"
Set-CellText $ws 8 3 "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-CellText $ws 8 4 "0.2868"
Set-CellText $ws 8 5 "  -0.05%  "

# Row 9
Set-CellText $ws 9 2 "OKB"
Set-CellText $ws 9 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-CellText $ws 9 4 "42.47"
Set-CellText $ws 9 5 "  -2.85%  "

# Row 10
Set-CellText $ws 10 4 "0.06468"
Set-CellText $ws 10 5 "  -1.40%  "

# Row 11
Set-CellText $ws 11 4 "20.88"
Set-CellText $ws 11 5 "  -3.53%  "

# Row 12
Set-CellText $ws 12 5 "  -1.94%  "

# Row 13
Set-CellText $ws 13 4 "1.868.66"
Set-CellText $ws 13 5 "  -0.04%  "

# Row 14
Set-CellText $ws 14 4 "95.02"
Set-CellText $ws 14 5 "  -1.11%  "

# Row 15
Set-CellText $ws 15 4 "0.7029"
Set-CellText $ws 15 5 "  +1.31%  "

# Row 16
Set-CellText $ws 16 4 "5.097"
Set-CellText $ws 16 5 "  +0.06%  "

# Row 17
Set-CellText $ws 17 4 "272.20"
Set-CellText $ws 17 5 "  +2.49%  "

# Row 18
Set-CellText $ws 18 4 "30.274.02"
Set-CellText $ws 18 5 "  +0.06%  "

# Row 19
Set-CellText $ws 19 4 "13.34"
Set-CellText $ws 19 5 "  -4.92%  "

# Row 20
Set-CellText $ws 20 4 "0.000007539"
Set-CellText $ws 20 5 "  -1.48%  "

# Row 21
Set-CellText $ws 21 5 "  -0.16%  "

# Row 22
Set-CellText $ws 22 4 "2.116.93"
Set-CellText $ws 22 5 "  +0.04%  "

# Row 23
Set-CellText $ws 23 4 "0.9997"
Set-CellText $ws 23 5 "  -0.17%  "

# Row 24
Set-CellText $ws 24 4 "5.208"
Set-CellText $ws 24 5 "  -0.44%  "

# Row 25
Set-CellText $ws 25 4 "6.127"
Set-CellText $ws 25 5 "  -0.98%  "

# Row 26
Set-CellText $ws 26 4 "9.312"
Set-CellText $ws 26 5 "  -0.63%  "

# Row 27
Set-CellText $ws 27 4 "165.31"
Set-CellText $ws 27 5 "  -1.15%  "

# Row 28
Set-CellText $ws 28 4 "18.87"
Set-CellText $ws 28 5 "  +0.51%  "

# Row 29
Set-CellText $ws 29 4 "1.903"
Set-CellText $ws 29 5 "  -1.87%  "

# Row 30
Set-CellText $ws 30 5 "  +1.66%  "

# Row 31
Set-CellText $ws 31 4 "0.09859"
Set-CellText $ws 31 5 "  -0.15%  "

# Row 32
Set-CellText $ws 32 4 "1.508"
Set-CellText $ws 32 5 "  +3.44%  "

# Row 33
Set-CellText $ws 33 5 "  -2.67%  "

# Row 34
Set-CellText $ws 34 4 "4.020"
Set-CellText $ws 34 5 "  -0.96%  "

# Row 35
Set-CellText $ws 35 4 "0.04725"
Set-CellText $ws 35 5 "  -0.62%  "

# Row 36
Set-CellText $ws 36 4 "1.121"
Set-CellText $ws 36 5 "  -0.98%  "

# Row 37
Set-CellText $ws 37 4 "0.6910"
Set-CellText $ws 37 5 "  -1.36%  "

# Row 38
Set-CellText $ws 38 4 "2.703"
Set-CellText $ws 38 5 "  -0.83%  "

# Row 39
Set-CellText $ws 39 4 "0.01843"
Set-CellText $ws 39 5 "  -1.41%  "

# Row 40
Set-CellText $ws 40 4 "2.742"
Set-CellText $ws 40 5 "  -1.89%  "

# Row 41
Set-CellText $ws 41 4 "6.305"
Set-CellText $ws 41 5 "  +1.62%  "

# Row 42
Set-CellText $ws 42 4 "70.08"
Set-CellText $ws 42 5 "  -3.33%  "

# Row 43
Set-CellText $ws 43 4 "0.9997"
Set-CellText $ws 43 5 "  -0.17%  "

# Row 44
Set-CellText $ws 44 4 "0.8394"
Set-CellText $ws 44 5 "  -0.36%  "

# Row 45
Set-CellText $ws 45 4 "1.896"
Set-CellText $ws 45 5 "  -2.28%  "

# Row 46
Set-CellText $ws 46 4 "102.03"
Set-CellText $ws 46 5 "  -0.34%  "

# Row 47
Set-CellText $ws 47 4 "0.4076"
Set-CellText $ws 47 5 "  -2.10%  "

# Row 48
Set-CellText $ws 48 4 "9.244"
Set-CellText $ws 48 5 "  +1.78%  "

# Row 49
Set-CellText $ws 49 4 "7.059"
Set-CellText $ws 49 5 "  -0.80%  "

# Row 50
Set-CellText $ws 50 4 "924.86"
Set-CellText $ws 50 5 "  -1.37%  "

# Row 51
Set-CellText $ws 51 4 "34.78"
Set-CellText $ws 51 5 "  +0.94%  "
